$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1240.53
$ws.Range("I15").Value = 1240.53
$ws.Range("K15").Value = 3721.59
$ws.Range("M15").Value = -3552.59
$ws.Range("H62").Value = 2108.3333
$ws.Range("J62").Value = 1749
$ws.Range("L62").Value = 1749
$ws.Range("N62").Value = -2997
$ws.Range("H65").Value = 2108.3333
$ws.Range("J65").Value = 1749
$ws.Range("L65").Value = 8745
$ws.Range("N65").Value = -14985
$ws.Range("I113").Value = 2890
$ws.Range("J113").Value = 3780.5
$ws.Range("K113").Value = 2890
$ws.Range("L113").Value = 3780.5
$ws.Range("M113").Value = 364
$ws.Range("N113").Value = -10288.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3832.5366
$ws.Range("J63").Value = 10244.25
$ws.Range("L63").Value = 10244.25
$ws.Range("N63").Value = -11616.25
$ws.Range("H66").Value = 3832.5366
$ws.Range("J66").Value = 10244.25
$ws.Range("L66").Value = 51221.25
$ws.Range("N66").Value = -58085.25
$ws.Range("H101").Value = 16379.8
$ws.Range("J101").Value = 16379.8
$ws.Range("L101").Value = 16379.8
$ws.Range("N101").Value = -22869.8
$ws.Range("H105").Value = 34932.668
$ws.Range("J105").Value = 34932.668
$ws.Range("L105").Value = 34932.668
$ws.Range("N105").Value = -41920.668
$ws.Range("H118").Value = 39650
$ws.Range("J118").Value = 39650
$ws.Range("L118").Value = 39650
$ws.Range("N118").Value = -42964
$ws.Range("H132").Value = 1347130
$ws.Range("I132").Value = 1940.9032
$ws.Range("J132").Value = 4037508.2
$ws.Range("K132").Value = 5822.7096
$ws.Range("L132").Value = 12112524.6
$ws.Range("M132").Value = -3292.7096
$ws.Range("N132").Value = -12117584.6

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2149.9473
$ws.Range("I94").Value = 2060.6428
$ws.Range("J94").Value = 2400
$ws.Range("K94").Value = 2060.6428
$ws.Range("L94").Value = 2400
$ws.Range("M94").Value = -1609.6428
$ws.Range("N94").Value = -3302
$ws.Range("H109").Value = 29356.666
$ws.Range("J109").Value = 29356.666
$ws.Range("L109").Value = 29356.666
$ws.Range("N109").Value = -32130.666
$ws.Range("H116").Value = 30742
$ws.Range("J116").Value = 30742
$ws.Range("L116").Value = 30742
$ws.Range("N116").Value = -39920
$ws.Range("H134").Value = 2466.0212
$ws.Range("I134").Value = 1395.2858
$ws.Range("J134").Value = 3330.8462
$ws.Range("K134").Value = 4185.857400000001
$ws.Range("L134").Value = 9992.5386
$ws.Range("M134").Value = -1650.857400000001
$ws.Range("N134").Value = -15062.5386

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1184.87
$ws.Range("I31").Value = 928.3333
$ws.Range("J31").Value = 2531.6875
$ws.Range("K31").Value = 928.3333
$ws.Range("L31").Value = 2531.6875
$ws.Range("M31").Value = -633.3333
$ws.Range("N31").Value = -3121.6875
$ws.Range("H34").Value = 1184.87
$ws.Range("I34").Value = 928.3333
$ws.Range("J34").Value = 2531.6875
$ws.Range("K34").Value = 928.3333
$ws.Range("L34").Value = 2531.6875
$ws.Range("M34").Value = -726.3333
$ws.Range("N34").Value = -2935.6875
$ws.Range("H58").Value = 4551.8823
$ws.Range("I58").Value = 3189.5881
$ws.Range("J58").Value = 5914.1763
$ws.Range("K58").Value = 3189.5881
$ws.Range("L58").Value = 5914.1763
$ws.Range("M58").Value = -2986.5881
$ws.Range("N58").Value = -6320.1763
$ws.Range("H74").Value = 20525.5
$ws.Range("J74").Value = 20525.5
$ws.Range("L74").Value = 20525.5
$ws.Range("N74").Value = -22273.5
$ws.Range("H77").Value = 20525.5
$ws.Range("J77").Value = 20525.5
$ws.Range("L77").Value = 61576.5
$ws.Range("N77").Value = -70312.5
$ws.Range("H107").Value = 1821.1428
$ws.Range("I107").Value = 452.5
$ws.Range("J107").Value = 2368.6
$ws.Range("K107").Value = 452.5
$ws.Range("L107").Value = 2368.6
$ws.Range("M107").Value = 1467.5
$ws.Range("N107").Value = -6208.6
$ws.Range("H136").Value = 4551.8823
$ws.Range("I136").Value = 3189.5881
$ws.Range("J136").Value = 5914.1763
$ws.Range("K136").Value = 9568.764299999999
$ws.Range("L136").Value = 17742.5289
$ws.Range("M136").Value = -7018.764299999999
$ws.Range("N136").Value = -22842.5289

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9371
$ws.Range("J3").Value = 12272.857
$ws.Range("L3").Value = 36818.571
$ws.Range("N3").Value = -37042.571
$ws.Range("H129").Value = 1900.125
$ws.Range("I129").Value = 728.3333
$ws.Range("J129").Value = 2603.2
$ws.Range("K129").Value = 2184.9999
$ws.Range("L129").Value = 7809.599999999999
$ws.Range("M129").Value = 2815.0001
$ws.Range("N129").Value = -17809.6
$ws.Range("H131").Value = 305668.28
$ws.Range("I131").Value = 458.42105
$ws.Range("J131").Value = 569258.5600000001
$ws.Range("K131").Value = 1375.26315
$ws.Range("L131").Value = 1707775.68
$ws.Range("M131").Value = 3664.73685
$ws.Range("N131").Value = -1717855.68
$ws.Range("H133").Value = 4238.36
$ws.Range("I133").Value = 2800
$ws.Range("J133").Value = 4298.2915
$ws.Range("K133").Value = 8400
$ws.Range("L133").Value = 12894.8745
$ws.Range("M133").Value = -3340
$ws.Range("N133").Value = -23014.8745
$ws.Range("H134").Value = 2870.6572
$ws.Range("I134").Value = 1395.3846
$ws.Range("J134").Value = 3742.4092
$ws.Range("K134").Value = 4186.1538
$ws.Range("L134").Value = 11227.2276
$ws.Range("M134").Value = 883.8462
$ws.Range("N134").Value = -21367.2276
$ws.Range("H136").Value = 4552.381
$ws.Range("I136").Value = 2060
$ws.Range("J136").Value = 6818.1816
$ws.Range("K136").Value = 6180
$ws.Range("L136").Value = 20454.5448
$ws.Range("M136").Value = -1080
$ws.Range("N136").Value = -30654.5448
$ws.Range("H138").Value = 2983.7144
$ws.Range("I138").Value = 1650.6154
$ws.Range("J138").Value = 5150
$ws.Range("K138").Value = 4951.8462
$ws.Range("L138").Value = 15450
$ws.Range("M138").Value = 188.1538
$ws.Range("N138").Value = -25730
$ws.Range("H139").Value = 2746.5405
$ws.Range("I139").Value = 501
$ws.Range("J139").Value = 2808.9167
$ws.Range("K139").Value = 1503
$ws.Range("L139").Value = 8426.750100000001
$ws.Range("M139").Value = 3637
$ws.Range("N139").Value = -18706.7501
$ws.Range("H140").Value = 1426.0526
$ws.Range("I140").Value = 1221.0714
$ws.Range("K140").Value = 3663.2142
$ws.Range("M140").Value = 1516.7858

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 15887.5
$ws.Range("J123").Value = 15887.5
$ws.Range("L123").Value = 15887.5
$ws.Range("N123").Value = -20787.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H68").Value = 1946.6666
$ws.Range("I68").Value = 1630.7693
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 1630.7693
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -881.7692999999999
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 1946.6666
$ws.Range("I71").Value = 1630.7693
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 8153.8465
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -4409.8465
$ws.Range("N71").Value = -27488
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H111").Value = 29546.75
$ws.Range("J111").Value = 29546.75
$ws.Range("L111").Value = 29546.75
$ws.Range("N111").Value = -37726.75

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 444.82352
$ws.Range("I107").Value = 267.07693
$ws.Range("J107").Value = 1022.5
$ws.Range("K107").Value = 801.2307900000001
$ws.Range("L107").Value = 3067.5
$ws.Range("M107").Value = 1118.76921
$ws.Range("N107").Value = -6907.5
$ws.Range("H121").Value = 36300
$ws.Range("J121").Value = 36300
$ws.Range("L121").Value = 36300
$ws.Range("N121").Value = -39794
$ws.Range("H132").Value = 1777.625
$ws.Range("I132").Value = 1512.0312
$ws.Range("K132").Value = 4536.0936
$ws.Range("M132").Value = -2006.0936
